$wb = $excel.ActiveWorkbook

# --- Content edits: rename underscore-separated transaction-type labels to
#     hyphen-separated ones (GOODS_RECEIPT -> GOODS-RECEIPT, etc.) ---

$wsInventory = $wb.Worksheets.Item("INVENTORY")
$wsInventory.Range("A3").Value = "GOODS-RECEIPT"
$wsInventory.Range("A4").Value = "GOODS-ISSUE"

$wsBusinessTransaction = $wb.Worksheets.Item("BUSINESS-TRANSACTION")
$wsBusinessTransaction.Range("A2").Value = "CUSTOMER-SALES"
$wsBusinessTransaction.Range("C3").Value = "CUSTOMER-SALES"

# --- View-state edits: selection cursor per sheet and active tab ---

$wsFishDispatch = $wb.Worksheets.Item("FISH-DISPATCH")
$wsFishDispatch.Range("C41").Select()

$wsInventory.Range("H15").Select()

$wsBusinessTransaction.Range("D21").Select()

# FISH-DISPATCH becomes the active (selected) tab/sheet.
$wsFishDispatch.Activate()
